# Daily attendance processing - 2026-01-24 08:38:55
# Rotate the "Recorded By" (column G) comma-separated list left by one
# entry (the first recorder moves to the end) for every data row that
# has more than one recorder listed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $value = $cell.Value()

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "
        if ($parts.Count -gt 1) {
            $rest = $parts[1..($parts.Count - 1)]
            $rotated = $rest + $parts[0]
            $cell.Value = [string]::Join(", ", $rotated)
        }
    }
}
